# Update NATMI LR-pair TPM values (columns E through T, rows 2-10)
# to reflect the new TPM-based calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> values for columns E..T (16 values per row)
$data = @{
    2 = @(2, 0.6666666666666666, 0.3851733333333334, 1.15552, 0.2330623765375699, 0.2330623765375699, 1, 0.3333333333333333, 0.06761166666666667, 0.202835, 0.04506879582203475, 0.04506879582203475, 0.02604221102222222, 0.2343798992, 0.01050384066196992, 0.01050384066196992)
    3 = @(2, 0.6666666666666666, 0.3851733333333334, 1.15552, 0.2330623765375699, 0.2330623765375699, 3, 1, 1.211564666666667, 3.634694, 0.8076085575052372, 0.8076085575052372, 0.4666624012088889, 4.199961610880001, 0.1882231697242493, 0.1882231697242493)
    4 = @(2, 0.6666666666666666, 0.3851733333333334, 1.15552, 0.2330623765375699, 0.2330623765375699, 3, 1, 0.2210116666666666, 0.6630349999999999, 0.1473226466727281, 0.1473226466727281, 0.08512780035555555, 0.7661502031999999, 0.03433536615135073, 0.03433536615135074)
    5 = @(3, 1, 1.032613, 3.097839, 0.6248180208657306, 0.6248180208657306, 1, 0.3333333333333333, 0.06761166666666667, 0.202835, 0.04506879582203475, 0.04506879582203475, 0.06981668595166667, 0.628350173565, 0.02815979580832546, 0.02815979580832546)
    6 = @(3, 1, 1.032613, 3.097839, 0.6248180208657306, 0.6248180208657306, 3, 1, 1.211564666666667, 3.634694, 0.8076085575052372, 0.8076085575052372, 1.251077425140667, 11.259696826266, 0.5046083805346498, 0.5046083805346498)
    7 = @(3, 1, 1.032613, 3.097839, 0.6248180208657306, 0.6248180208657306, 3, 1, 0.2210116666666666, 0.6630349999999999, 0.1473226466727281, 0.1473226466727281, 0.2282195201516666, 2.053975681365, 0.09204984452275529, 0.09204984452275529)
    8 = @(2, 0.6666666666666666, 0.2348756666666667, 0.704627, 0.1421196025966995, 0.1421196025966995, 1, 0.3333333333333333, 0.06761166666666667, 0.202835, 0.04506879582203475, 0.04506879582203475, 0.01588033528277778, 0.142923017545, 0.006405159351739372, 0.006405159351739371)
    9 = @(2, 0.6666666666666666, 0.2348756666666667, 0.704627, 0.1421196025966995, 0.1421196025966995, 3, 1, 1.211564666666667, 3.634694, 0.8076085575052372, 0.8076085575052372, 0.2845670587931111, 2.561103529138, 0.1147770072463381, 0.1147770072463381)
    10 = @(2, 0.6666666666666666, 0.2348756666666667, 0.704627, 0.1421196025966995, 0.1421196025966995, 3, 1, 0.2210116666666666, 0.6630349999999999, 0.1473226466727281, 0.1473226466727281, 0.05191026254944444, 0.4671923629449999, 0.0209374359986221, 0.0209374359986221)
}

# Columns E (5) through T (20)
$startCol = 5

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $startCol + $i).Value = $vals[$i]
    }
}
